$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3399353333333333
$ws.Range("H2").Value = 1.019806
$ws.Range("I2").Value = 0.09929991924017606
$ws.Range("J2").Value = 0.09929991924017606
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 0.9259844145588888
$ws.Range("R2").Value = 8.33385973103
$ws.Range("S2").Value = 0.004603048227294915
$ws.Range("T2").Value = 0.004603048227294915
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3399353333333333
$ws.Range("H3").Value = 1.019806
$ws.Range("I3").Value = 0.09929991924017606
$ws.Range("J3").Value = 0.09929991924017606
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 13.80838762588356
$ws.Range("R3").Value = 124.275488632952
$ws.Range("S3").Value = 0.06864119220991721
$ws.Range("T3").Value = 0.06864119220991721
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3399353333333333
$ws.Range("H4").Value = 1.019806
$ws.Range("I4").Value = 0.09929991924017606
$ws.Range("J4").Value = 0.09929991924017606
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 5.241559786236667
$ws.Range("R4").Value = 47.17403807613
$ws.Range("S4").Value = 0.02605567880296393
$ws.Range("T4").Value = 0.02605567880296393
$ws.Range("I5").Value = 0.4094685684206303
$ws.Range("J5").Value = 0.4094685684206303
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 3.81834663623611
$ws.Range("R5").Value = 34.36511972612499
$ws.Range("S5").Value = 0.01898091743098811
$ws.Range("T5").Value = 0.01898091743098811
$ws.Range("I6").Value = 0.4094685684206303
$ws.Range("J6").Value = 0.4094685684206303
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("S6").Value = 0.2830456552628138
$ws.Range("T6").Value = 0.2830456552628138
$ws.Range("I7").Value = 0.4094685684206303
$ws.Range("J7").Value = 0.4094685684206303
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 21.61385425470833
$ws.Range("R7").Value = 194.524688292375
$ws.Range("S7").Value = 0.1074419957268284
$ws.Range("T7").Value = 0.1074419957268284
$ws.Range("G8").Value = 1.681642333333333
$ws.Range("H8").Value = 5.044927
$ws.Range("I8").Value = 0.4912315123391937
$ws.Range("J8").Value = 0.4912315123391937
$ws.Range("M8").Value = 2.724001666666667
$ws.Range("N8").Value = 8.172005
$ws.Range("O8").Value = 0.04635500474236593
$ws.Range("P8").Value = 0.04635500474236593
$ws.Range("Q8").Value = 4.580796518737222
$ws.Range("R8").Value = 41.227168668635
$ws.Range("S8").Value = 0.02277103908408291
$ws.Range("T8").Value = 0.02277103908408291
$ws.Range("G9").Value = 1.681642333333333
$ws.Range("H9").Value = 5.044927
$ws.Range("I9").Value = 0.4912315123391937
$ws.Range("J9").Value = 0.4912315123391937
$ws.Range("O9").Value = 0.6912512390256352
$ws.Range("P9").Value = 0.6912512390256351
$ws.Range("Q9").Value = 68.3093721357649
$ws.Range("R9").Value = 614.7843492218841
$ws.Range("S9").Value = 0.3395643915529042
$ws.Range("T9").Value = 0.3395643915529042
$ws.Range("G10").Value = 1.681642333333333
$ws.Range("H10").Value = 5.044927
$ws.Range("I10").Value = 0.4912315123391937
$ws.Range("J10").Value = 0.4912315123391937
$ws.Range("M10").Value = 15.419285
$ws.Range("N10").Value = 46.257855
$ws.Range("O10").Value = 0.2623937562319988
$ws.Range("P10").Value = 0.2623937562319988
$ws.Range("Q10").Value = 25.92972240573167
$ws.Range("R10").Value = 233.367501651585
$ws.Range("S10").Value = 0.1288960817022065
$ws.Range("T10").Value = 0.1288960817022065
